$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared-strings table had a new string "LOT2038 - Tecnologia de Bebidas
# (Indicacao de Conjunto)" inserted immediately before the existing
# "LOT2028 - Tecnologia de Processos Fermentativos (Requisito fraco)" string.
# Since the worksheet's cell references did not change, the practical effect
# is that the two requirement rows swap their displayed text:
#   Row 23 (previously LOT2028 ...) now shows the LOT2038 ... text
#   Row 24 (previously LOT2038 ...) now shows the LOT2028 ... text

$lot2028 = "LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)`n"
$lot2038 = "LOT2038 -  Tecnologia de Bebidas  (Indicação de Conjunto)`n"

$ws.Range("B23").Value = $lot2038
$ws.Range("C23").Value = $lot2038

$ws.Range("B24").Value = $lot2028
$ws.Range("C24").Value = $lot2028
